$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 69, mirroring the constant columns of row 68, then fill shifted values ---
$ws.Range("A69").Value = 10
$ws.Range("B69").Value = "Vega Modelo de Temuco"
$ws.Range("C69").Value = "La Araucanía"
$ws.Range("D68").Copy()
$ws.Range("D69").PasteSpecial(-4122)
$ws.Range("D69").Value = 44519
$ws.Range("E69").Value = 9
$ws.Range("F69").Value = 100112026
$ws.Range("G69").Value = "Haba"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 55
$ws.Range("K69").Value = 8000
$ws.Range("L69").Value = 9000
$ws.Range("M69").Value = 8455
$ws.Range("N69").Value = "$/saco 25 kilos"
$ws.Range("O69").Value = "Región del Maule"
$ws.Range("P69").Value = 338
$ws.Range("Q69").Value = 25
$ws.Range("R69").Value = "Hortaliza"

# --- Shift rows 33-68 down by one (values that were in row N now sit in row N+1); row 33 gets the new weekly entry ---
$ws.Range("D33").Value = 44740
$ws.Range("J33").Value = 40
$ws.Range("K33").Value = 20000
$ws.Range("L33").Value = 20000
$ws.Range("M33").Value = 20000
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 800

$ws.Range("D34").Value = 44427
$ws.Range("J34").Value = 30
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = 15000
$ws.Range("O34").Value = "Provincia de Limarí"
$ws.Range("P34").Value = 600

$ws.Range("D35").Value = 44468
$ws.Range("J35").Value = 30
$ws.Range("K35").Value = 14000
$ws.Range("L35").Value = 14000
$ws.Range("M35").Value = 14000
$ws.Range("O35").Value = "Provincia de Limarí"
$ws.Range("P35").Value = 560

$ws.Range("D36").Value = 44434
$ws.Range("J36").Value = 50
$ws.Range("K36").Value = 15000
$ws.Range("L36").Value = 15000
$ws.Range("M36").Value = 15000
$ws.Range("O36").Value = "Provincia de Limarí"
$ws.Range("P36").Value = 600

$ws.Range("D37").Value = 44448
$ws.Range("J37").Value = 35
$ws.Range("K37").Value = 15000
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = 15000
$ws.Range("O37").Value = "Provincia de Limarí"
$ws.Range("P37").Value = 600

$ws.Range("D38").Value = 44186
$ws.Range("J38").Value = 30
$ws.Range("K38").Value = 14000
$ws.Range("L38").Value = 14000
$ws.Range("M38").Value = 14000
$ws.Range("O38").Value = "Región de La Araucanía"
$ws.Range("P38").Value = 560

$ws.Range("D39").Value = 44515
$ws.Range("J39").Value = 155
$ws.Range("K39").Value = 9000
$ws.Range("L39").Value = 9000
$ws.Range("M39").Value = 9000
$ws.Range("O39").Value = "Región del Maule"
$ws.Range("P39").Value = 360

$ws.Range("D40").Value = 44160
$ws.Range("J40").Value = 30
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = 8000
$ws.Range("O40").Value = "Región de La Araucanía"
$ws.Range("P40").Value = 320

$ws.Range("D41").Value = 44487
$ws.Range("J41").Value = 110
$ws.Range("K41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = 10000
$ws.Range("O41").Value = "Región Metropolitana"
$ws.Range("P41").Value = 400

$ws.Range("D42").Value = 44496
$ws.Range("J42").Value = 50
$ws.Range("K42").Value = 8000
$ws.Range("L42").Value = 9000
$ws.Range("M42").Value = 8600
$ws.Range("O42").Value = "Región Metropolitana"
$ws.Range("P42").Value = 344

$ws.Range("D43").Value = 44526
$ws.Range("J43").Value = 20
$ws.Range("K43").Value = 8000
$ws.Range("L43").Value = 8000
$ws.Range("M43").Value = 8000
$ws.Range("O43").Value = "Región de La Araucanía"
$ws.Range("P43").Value = 320

$ws.Range("D44").Value = 44490
$ws.Range("J44").Value = 65
$ws.Range("K44").Value = 9000
$ws.Range("L44").Value = 9000
$ws.Range("M44").Value = 9000
$ws.Range("O44").Value = "Región Metropolitana"
$ws.Range("P44").Value = 360

$ws.Range("D45").Value = 44371
$ws.Range("J45").Value = 40
$ws.Range("K45").Value = 15000
$ws.Range("L45").Value = 15000
$ws.Range("M45").Value = 15000
$ws.Range("O45").Value = "Provincia de Limarí"
$ws.Range("P45").Value = 600

$ws.Range("D46").Value = 44210
$ws.Range("J46").Value = 110
$ws.Range("K46").Value = 16000
$ws.Range("L46").Value = 16000
$ws.Range("M46").Value = 16000
$ws.Range("O46").Value = "Región de La Araucanía"
$ws.Range("P46").Value = 640

$ws.Range("D47").Value = 44484
$ws.Range("J47").Value = 30
$ws.Range("K47").Value = 9000
$ws.Range("L47").Value = 9000
$ws.Range("M47").Value = 9000
$ws.Range("O47").Value = "Provincia de Limarí"
$ws.Range("P47").Value = 360

$ws.Range("D48").Value = 44455
$ws.Range("J48").Value = 10
$ws.Range("K48").Value = 13000
$ws.Range("L48").Value = 13000
$ws.Range("M48").Value = 13000
$ws.Range("O48").Value = "Provincia de Limarí"
$ws.Range("P48").Value = 520

$ws.Range("D49").Value = 44516
$ws.Range("J49").Value = 85
$ws.Range("K49").Value = 9000
$ws.Range("L49").Value = 9000
$ws.Range("M49").Value = 9000
$ws.Range("O49").Value = "Región del Maule"
$ws.Range("P49").Value = 360

$ws.Range("D50").Value = 44491
$ws.Range("J50").Value = 55
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 10000
$ws.Range("M50").Value = 10000
$ws.Range("O50").Value = "Región Metropolitana"
$ws.Range("P50").Value = 400

$ws.Range("D51").Value = 44356
$ws.Range("J51").Value = 30
$ws.Range("K51").Value = 14000
$ws.Range("L51").Value = 14000
$ws.Range("M51").Value = 14000
$ws.Range("O51").Value = "Provincia de Limarí"
$ws.Range("P51").Value = 560

$ws.Range("D52").Value = 44720
$ws.Range("J52").Value = 40
$ws.Range("K52").Value = 25000
$ws.Range("L52").Value = 25000
$ws.Range("M52").Value = 25000
$ws.Range("O52").Value = "Provincia de Limarí"
$ws.Range("P52").Value = 1000

$ws.Range("D53").Value = 44452
$ws.Range("J53").Value = 50
$ws.Range("K53").Value = 13000
$ws.Range("L53").Value = 13000
$ws.Range("M53").Value = 13000
$ws.Range("O53").Value = "Provincia de Limarí"
$ws.Range("P53").Value = 520

$ws.Range("D54").Value = 44469
$ws.Range("J54").Value = 90
$ws.Range("K54").Value = 14000
$ws.Range("L54").Value = 14000
$ws.Range("M54").Value = 14000
$ws.Range("O54").Value = "Provincia de Limarí"
$ws.Range("P54").Value = 560

$ws.Range("D55").Value = 44435
$ws.Range("J55").Value = 50
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = 15000
$ws.Range("O55").Value = "Provincia de Limarí"
$ws.Range("P55").Value = 600

$ws.Range("D56").Value = 44357
$ws.Range("J56").Value = 80
$ws.Range("K56").Value = 14000
$ws.Range("L56").Value = 14000
$ws.Range("M56").Value = 14000
$ws.Range("O56").Value = "Provincia de Limarí"
$ws.Range("P56").Value = 560

$ws.Range("D57").Value = 44476
$ws.Range("J57").Value = 50
$ws.Range("K57").Value = 10000
$ws.Range("L57").Value = 10000
$ws.Range("M57").Value = 10000
$ws.Range("O57").Value = "Provincia de Limarí"
$ws.Range("P57").Value = 400

$ws.Range("D58").Value = 44482
$ws.Range("J58").Value = 40
$ws.Range("K58").Value = 10000
$ws.Range("L58").Value = 10000
$ws.Range("M58").Value = 10000
$ws.Range("O58").Value = "Provincia de Limarí"
$ws.Range("P58").Value = 400

$ws.Range("D59").Value = 44474
$ws.Range("J59").Value = 30
$ws.Range("K59").Value = 10000
$ws.Range("L59").Value = 10000
$ws.Range("M59").Value = 10000
$ws.Range("O59").Value = "Provincia de Limarí"
$ws.Range("P59").Value = 400

$ws.Range("D60").Value = 44162
$ws.Range("J60").Value = 260
$ws.Range("K60").Value = 7000
$ws.Range("L60").Value = 8000
$ws.Range("M60").Value = 7462
$ws.Range("O60").Value = "Región de La Araucanía"
$ws.Range("P60").Value = 298

$ws.Range("D61").Value = 44494
$ws.Range("J61").Value = 100
$ws.Range("K61").Value = 9000
$ws.Range("L61").Value = 9000
$ws.Range("M61").Value = 9000
$ws.Range("O61").Value = "Región Metropolitana"
$ws.Range("P61").Value = 360

$ws.Range("D62").Value = 44553
$ws.Range("J62").Value = 175
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = 10971
$ws.Range("O62").Value = "Región de La Araucanía"
$ws.Range("P62").Value = 439

$ws.Range("D63").Value = 44473
$ws.Range("J63").Value = 80
$ws.Range("K63").Value = 10000
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = 10000
$ws.Range("O63").Value = "Provincia de Limarí"
$ws.Range("P63").Value = 400

$ws.Range("D64").Value = 44159
$ws.Range("J64").Value = 50
$ws.Range("K64").Value = 8000
$ws.Range("L64").Value = 8000
$ws.Range("M64").Value = 8000
$ws.Range("O64").Value = "Región de La Araucanía"
$ws.Range("P64").Value = 320

$ws.Range("D65").Value = 44159
$ws.Range("J65").Value = 80
$ws.Range("K65").Value = 8000
$ws.Range("L65").Value = 8000
$ws.Range("M65").Value = 8000
$ws.Range("O65").Value = "Región del Maule"
$ws.Range("P65").Value = 320

$ws.Range("D66").Value = 44518
$ws.Range("J66").Value = 125
$ws.Range("K66").Value = 8000
$ws.Range("L66").Value = 8000
$ws.Range("M66").Value = 8000
$ws.Range("O66").Value = "Región del Maule"
$ws.Range("P66").Value = 320

$ws.Range("D67").Value = 44505
$ws.Range("J67").Value = 120
$ws.Range("K67").Value = 7000
$ws.Range("L67").Value = 7000
$ws.Range("M67").Value = 7000
$ws.Range("O67").Value = "Provincia de Limarí"
$ws.Range("P67").Value = 280

$ws.Range("D68").Value = 44483
$ws.Range("J68").Value = 80
$ws.Range("K68").Value = 9000
$ws.Range("L68").Value = 9000
$ws.Range("M68").Value = 9000
$ws.Range("O68").Value = "Provincia de Limarí"
$ws.Range("P68").Value = 360
